# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 11:05"

# --- 2. Alemania (row 11): Casos activos / Recuperados corrected ---
$ws.Range("D11").Value = 159900
$ws.Range("E11").Value = 11461

# --- 3. Banglades: move up in the ranking with fresh stats.
#    Currently Banglades sits below Suiza (row 30); it needs to move to just
#    above Singapur (currently row 28) with updated numbers. Insert a new
#    row there, fill it in, then remove the old Banglades row (now pushed
#    down to row 31).
$ws.Rows("28:28").Insert()
$ws.Range("A28").Value = "Banglades"
$ws.Range("B28").Value = 32078
$ws.Range("C28").Value = 1873
$ws.Range("D28").Value = 6486
$ws.Range("E28").Value = 25140
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = 452

# The old Banglades row has been pushed down to row 31; delete it.
$ws.Rows("31:31").Delete()

# --- 4. Malasia (row 60): updated stats ---
$ws.Range("B60").Value = 7185
$ws.Range("C60").Value = 48
$ws.Range("D60").Value = 5912
$ws.Range("E60").Value = 1158

# --- 5. Consejo Danes para los Refugiados (row 87): updated stats ---
$ws.Range("B87").Value = 2025
$ws.Range("C87").Value = 80
$ws.Range("E87").Value = 1650

# --- 6. Kenia (row 102): updated stats ---
$ws.Range("B102").Value = 1192
$ws.Range("C102").Value = 31
$ws.Range("E102").Value = 762
